# DPLKKPS013-001 - add "SIDEBAR_SUBMENU_SUBMENU" column
# (Setup Kelengkapan Kepesertaan) between SIDEBAR_SUBMENU and KODE_JENIS_BIAYA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (KODE_JENIS_BIAYA), shifting the
# existing L/M/N columns (SIDEBAR_SUBMENU, KODE_JENIS_BIAYA, NAMA_JENIS_BIAYA
# values) one slot to the right. Excel carries the formatting of the column
# to the left (K) into the new column automatically.
$ws.Columns("L:L").Insert()

# Header (row 1) + data (row 2) for the newly inserted column.
$ws.Range("L1").Value = "SIDEBAR_SUBMENU_SUBMENU"
$ws.Range("L2").Value = "'Setup Kelengkapan Kepesertaan"

# Match the column widths to the new content/layout.
$ws.Columns("K").ColumnWidth = 12.333333333333334
$ws.Columns("L").ColumnWidth = 25.666666666666668
$ws.Columns("M").ColumnWidth = 24.166666666666668

# Update the active selection to reflect where the author left off.
$ws.Range("M2").Select()
